# Insert a new row at row 38, shifting existing rows 38..126 down to 39..127
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new data record
$ws.Cells.Item(38, 1).Value = 6
$ws.Cells.Item(38, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(38, 3).Value = "Metropolitana"
$ws.Cells.Item(38, 4).Value = 45054
$ws.Cells.Item(38, 5).Value = 13
$ws.Cells.Item(38, 6).Value = 100114007
$ws.Cells.Item(38, 7).Value = "Jengibre"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 290
$ws.Cells.Item(38, 11).Value = 15000
$ws.Cells.Item(38, 12).Value = 16000
$ws.Cells.Item(38, 13).Value = 15586
$ws.Cells.Item(38, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(38, 15).Value = "Perú"
$ws.Cells.Item(38, 16).Value = 1199
$ws.Cells.Item(38, 17).Value = 13
$ws.Cells.Item(38, 18).Value = "Hortaliza"
